# Penalty Reward System (unfinished) - weekly data rolled forward by one week
# and a couple of requested-quantity figures were corrected.

$wb = $excel.ActiveWorkbook

# --- "Weekly Quantity" sheet -------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# The oldest week (row 3: 2023-07-09) dropped off the report; deleting the
# entire row shifts every later week up by one and shrinks the used range
# from A1:B56 to A1:B55 automatically.
$wsWeekly.Rows.Item(3).Delete()

# Two of the requested-quantity figures (now at rows 26 and 27 after the
# shift) were revised.
$wsWeekly.Range("B26").Value = 420
$wsWeekly.Range("B27").Value = 20

# --- "Monthly Trend" sheet ------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Range("B3").Value = 320
$wsMonthly.Range("B11").Value = 1030
